# allDiagram.docx edit:
#  - Sirius provider options now use double quotes ("...") instead of
#    single quotes ('...') to enclose option values (AQL expression
#    support instead of a static String).
#  - The footer DATE field is collapsed from a begin/instrText/separate/
#    end run sequence into a <w:fldSimple> and its cached result is
#    refreshed (22/06/2016 -> 27/06/2016).
#  - The stray "_GoBack" bookmark that Word had left in the middle of
#    the first field's instrText runs is moved to the (until now empty)
#    paragraph that follows the second field.

$d = $word.ActiveDocument

# Helper: replace only the first occurrence of $old found at/after
# $script:pos in $script:xml, then advance $script:pos past the
# replacement so subsequent calls keep moving forward through the
# document instead of re-matching earlier (identical) text.
$script:xml = $d.Content.WordOpenXML
$script:pos = 0

function Replace-First([string]$old, [string]$new) {
    $idx = $script:xml.IndexOf($old, $script:pos)
    if ($idx -lt 0) {
        throw "pattern not found: $old"
    }
    $script:xml = $script:xml.Substring(0, $idx) + $new + $script:xml.Substring($idx + $old.Length)
    $script:pos = $idx + $new.Length
}

$runOpen = '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>'
$runClose = '</w:r>'

# ---------------------------------------------------------------
# Field 1: diagramProvider:'org.obeonetwork.m2doc.sirius.SiriusDiagramByRepresentationAndEObjectProvider' width:'200' height:'200' rootObject:'db.schemas->first()' diagramDescriptionName:'Schema Diagram'
# ---------------------------------------------------------------
Replace-First "<w:instrText>:'</w:instrText>" '<w:instrText>:"</w:instrText>'
Replace-First "<w:instrText>'</w:instrText>" '<w:instrText>"</w:instrText>'
Replace-First "<w:instrText>:'200' height:'200' rootObject:'</w:instrText>" '<w:instrText>:"200" height:"200" rootObject:"</w:instrText>'
Replace-First "<w:instrText>'</w:instrText>" '<w:instrText>"</w:instrText>'
Replace-First "<w:instrText>:'</w:instrText>" '<w:instrText>:"</w:instrText>'

# The "_GoBack" bookmark that used to sit between "Schema Diagram" and
# the closing quote run is removed here (it is re-created later, after
# the second field, on what used to be an empty paragraph).
Replace-First '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' ''
Replace-First "<w:instrText>'</w:instrText>" '<w:instrText>"</w:instrText>'

# ---------------------------------------------------------------
# Field 2: diagramProvider:'org.obeonetwork.m2doc.sirius.SiriusDiagramByTitleProvider' width:'200' height:'200' title:'RF Schema2'
# ---------------------------------------------------------------

# ":'"  ->  ":"  followed by a brand new run containing just the opening quote
Replace-First "<w:instrText>:'</w:instrText></w:r>" ("<w:instrText>:</w:instrText></w:r>" + $runOpen + '<w:instrText>"</w:instrText>' + $runClose)

Replace-First "<w:instrText>'</w:instrText>" '<w:instrText>"</w:instrText>'

# ":'200' height:'200' " (trailing space kept) -> ':"200" height:"200"' plus a
# new run holding just the trailing space
Replace-First '<w:instrText xml:space="preserve">:''200'' height:''200'' </w:instrText></w:r>' ('<w:instrText>:"200" height:"200"</w:instrText></w:r>' + $runOpen + '<w:instrText xml:space="preserve"> </w:instrText>' + $runClose)

# "'RF Schema" -> new run with just the opening quote, then a new run with "RF Schema"
Replace-First "<w:instrText>'RF Schema</w:instrText></w:r>" ('<w:instrText>"</w:instrText></w:r>' + $runOpen + '<w:instrText>RF Schema</w:instrText>' + $runClose)

Replace-First "<w:instrText>'</w:instrText>" '<w:instrText>"</w:instrText>'

# ---------------------------------------------------------------
# Trailing empty paragraph (after field 2) now carries the "_GoBack" bookmark
# ---------------------------------------------------------------
$oldTailParaOpen = ' w:rsidR="00577C6F" w:rsidRPr="00474E78" w:rsidRDefault="00577C6F" w:rsidP="00474E78"/>'
$newTailPara = ' w:rsidR="00577C6F" w:rsidRPr="00474E78" w:rsidRDefault="00577C6F" w:rsidP="00474E78"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
Replace-First $oldTailParaOpen $newTailPara

# ---------------------------------------------------------------
# Footer: collapse the DATE field's begin/instrText/separate/t/end runs
# into a <w:fldSimple> and refresh its cached result.
# ---------------------------------------------------------------
$oldDateField = '<w:r><w:fldChar w:fldCharType="begin"/></w:r>' `
    + '<w:r><w:instrText xml:space="preserve"> DATE   \* MERGEFORMAT </w:instrText></w:r>' `
    + '<w:r><w:fldChar w:fldCharType="separate"/></w:r>' `
    + '<w:r><w:rPr><w:noProof/></w:rPr><w:t>22/06/2016</w:t></w:r>' `
    + '<w:r><w:rPr><w:noProof/></w:rPr><w:fldChar w:fldCharType="end"/></w:r>'
$newDateField = '<w:fldSimple w:instr=" DATE   \* MERGEFORMAT ">' `
    + '<w:r><w:rPr><w:noProof/></w:rPr><w:t>27/06/2016</w:t></w:r>' `
    + '</w:fldSimple>'
$script:pos = 0
Replace-First $oldDateField $newDateField

$d.Content.WordOpenXML = $script:xml

Write-Output "done"
